$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Mark the empty "picture goes here" paragraph (the big cell under the
#    "СХЕМА ОБСЛЕДОВАННОГО УЧАСТКА" heading) with bookmarks "Image" and
#    "_GoBack" (nested: Image outer, _GoBack inner), so a picture can later
#    be inserted at that bookmark.
# ---------------------------------------------------------------------------
$rngHeading = $d.Content
$rngHeading.Find.Execute("СХЕМА ОБСЛЕДОВАННОГО УЧАСТКА", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingStart = $rngHeading.Start

$t = $d.Tables(1)
$headingRow = -1
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  $cell = $t.Cell($r, 1)
  if (($headingStart -ge $cell.Range.Start) -and ($headingStart -lt $cell.Range.End)) {
    $headingRow = $r
    break
  }
}

$imageCell = $t.Cell($headingRow + 1, 1)
$imagePara = $imageCell.Range.Paragraphs(1)
$imageTarget = $imagePara.Range

$bookmarkXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="006E0D42" w:rsidRPr="006E0D42" w:rsidRDefault="006E0D42" w:rsidP="002660E9"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:u w:val="single"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="Image"/><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:bookmarkEnd w:id="1"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$imageTarget.InsertXML($bookmarkXml)

# ---------------------------------------------------------------------------
# 2) Merge the two runs that spell out "{ispolnitel" + "}" (previously split
#    apart by the old "_GoBack" bookmark) into a single run with text
#    "{ispolnitel}", now that "_GoBack" has moved to the image placeholder.
# ---------------------------------------------------------------------------
$rngIspolnitel = $d.Content
$rngIspolnitel.Find.Execute("{ispolnitel}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ispolnitelPara = $rngIspolnitel.Paragraphs(1).Range

$mergedRunXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005A11D6" w:rsidRPr="002660E9" w:rsidRDefault="002660E9" w:rsidP="005855CE"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="32"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>{ispolnitel}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ispolnitelPara.InsertXML($mergedRunXml)

Write-Output "done"
